$d = $word.ActiveDocument

# --- Title line: date update ---
$d.Content.Find.Execute("2024-07-02 Tuesday", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "2024-07-03 Wednesday", 2) | Out-Null

# --- Table cell updates (row, col are 1-indexed Word table coordinates) ---
$tbl = $d.Tables(1)

function Set-Cell($row, $col, $text) {
    $cell = $tbl.Cell($row, $col)
    $cell.Range.Text = $text
}

# Data row 1 (Word table row 1)
Set-Cell 1 1 "71÷6=11, 5"
Set-Cell 1 2 "66÷2=33, 0"
Set-Cell 1 3 "52÷4=13, 0"
Set-Cell 1 4 "16÷2=8, 0"
Set-Cell 1 5 "52÷3=17, 1"

# Data row 2 (Word table row 5)
Set-Cell 5 1 "63÷9=7, 0"
Set-Cell 5 2 "20÷9=2, 2"
Set-Cell 5 3 "95÷6=15, 5"
Set-Cell 5 4 "35÷2=17, 1"
Set-Cell 5 5 "30÷9=3, 3"

# Data row 3 (Word table row 9)
Set-Cell 9 1 "68÷4=17, 0"
Set-Cell 9 2 "65÷7=9, 2"
Set-Cell 9 3 "69÷3=23, 0"
Set-Cell 9 4 "38÷2=19, 0"
Set-Cell 9 5 "43÷7=6, 1"

# Data row 4 (Word table row 13)
Set-Cell 13 1 "51÷2=25, 1"
Set-Cell 13 2 "89÷7=12, 5"
Set-Cell 13 3 "86÷2=43, 0"
Set-Cell 13 4 "88÷7=12, 4"
Set-Cell 13 5 "38÷7=5, 3"

# Data row 5 (Word table row 17)
Set-Cell 17 1 "78÷8=9, 6"
Set-Cell 17 2 "42÷9=4, 6"
Set-Cell 17 3 "68÷4=17, 0"
Set-Cell 17 4 "23÷6=3, 5"
Set-Cell 17 5 "95÷7=13, 4"
